$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column B, rows 2 through 45 (rows 46-97 remain 0, unchanged)
$bValues = @(2068,1511,1876,2019,1991,1963,1920,1875,1892,1835,1735,1676,1651,1647,1603,1539,1429,1405,1328,1279,1220,1123,1001,977,993,901,798,716,685,652,600,520,429,406,394,356,299,250,219,223,204,182,165,150)

# Column A (Timestamp) for rows 2-97 is shifted forward by exactly one day.
# Column B (Actual Production) gets new values for rows 2-45; rows 46-97 keep their existing 0 values.
for ($r = 2; $r -le 97; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.Value2 = $aCell.Value2 + 1

    if ($r -le 45) {
        $bCell = $ws.Cells.Item($r, 2)
        $bCell.Value2 = $bValues[$r - 2]
    }
}
